Write-Output "Direct: Não existente"
Write-Output "Descricao: Descrição de um produto"
